# Append a new data row (row 70) to each of the four worksheets, mirroring
# the layout/format of the existing last row (row 69). This corresponds to
# new log entries being appended to each "LIFTER" sheet.

$wb = $excel.ActiveWorkbook

$rowsToAdd = @(
    @{
        Sheet = "ROW35-FE-LIFTER"
        A = 45759.39265232639
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x66"
        E = "0xd"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 358
        I = 13
    },
    @{
        Sheet = "ROW35-MID-LIFTER"
        A = 45759.24753052084
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x66"
        E = "0xe"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 358
        I = 14
    },
    @{
        Sheet = "ROW02-FE-LIFTER"
        A = 45759.38840091435
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x01,0x66"
        E = "0x3"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 358
        I = 3
    },
    @{
        Sheet = "ROW02-MID-LIFTER"
        A = 45759.44550739583
        B = "0x01,0x90"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x01,0x66"
        E = "0x3"
        F = 400
        G = [double]"9.85046333984776e+23"
        H = 358
        I = 3
    }
)

foreach ($entry in $rowsToAdd) {
    $ws = $wb.Worksheets.Item($entry.Sheet)
    $newRow = $ws.Cells.Item(69, 1).EntireRow.Row + 1

    # Copy formatting from the row above (row 69) down into the new row,
    # then overwrite the values so the new row matches row 69's look & feel.
    $ws.Range("A69:I69").Copy() | Out-Null
    $ws.Range("A" + $newRow + ":I" + $newRow).PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($newRow, 1).Value = $entry.A
    $ws.Cells.Item($newRow, 2).Value = $entry.B
    $ws.Cells.Item($newRow, 3).Value = $entry.C
    $ws.Cells.Item($newRow, 4).Value = $entry.D
    $ws.Cells.Item($newRow, 5).Value = $entry.E
    $ws.Cells.Item($newRow, 6).Value = $entry.F
    $ws.Cells.Item($newRow, 7).Value = $entry.G
    $ws.Cells.Item($newRow, 8).Value = $entry.H
    $ws.Cells.Item($newRow, 9).Value = $entry.I
}

$excel.CutCopyMode = 0
